$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 112 previously carried the page-number "NA" in column C; the script
# re-ran and that value moved off of this row (cleared), matching every
# other "Rien ne nous concerne aujourd'hui !" row that has no page number.
$ws.Range("C112").Value = ""

# The new run appended a fresh row (113) for 2025-05-19 with the same
# "nothing relevant today" term and the "NA" page marker that used to sit
# on row 112.
$ws.Range("A113").NumberFormat = "@"
$ws.Range("A113").Value = "2025-05-19"
$ws.Range("A113").Style = "Normal"

$ws.Range("B113").Value = "Rien ne nous concerne aujourd'hui !"
$ws.Range("C113").Value = "NA"
$ws.Range("D113").Value = 1
